# Committed Corporate Customer excel file
# Update the Transaction Number value on Sheet0 and restore the cell
# selection left behind by the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet0")

# The transaction number in A2 was replaced with a new value.
$ws.Range("A2").Value = "FT2318504W72TY0Y"

# Leave the same cell selected as was left selected when the file was saved.
$ws.Range("M9").Select()
